# Update the cryptocurrency price/volume table to reflect the latest
# scrape (GitHub Actions scheduled run). Column D ("Price") holds values
# that look numeric (e.g. "1.015", "29.562.98") but must stay as literal
# text, exactly like the source data - so NumberFormat is forced to Text
# ("@") before assignment to stop Excel's automatic number coercion.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.562.98'
$ws.Range('E2').Value = '  -2.67%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.012.18'
$ws.Range('E3').Value = '  -4.62%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.015'
$ws.Range('E4').Value = '  +0.81%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.97'
$ws.Range('E5').Value = '  -3.86%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.012'
$ws.Range('E6').Value = '  +0.67%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5043'
$ws.Range('E7').Value = '  -3.65%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4258'
$ws.Range('E8').Value = '  -4.19%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.03'
$ws.Range('E9').Value = '  -0.95%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09237'
$ws.Range('E10').Value = '  -2.39%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.127'
$ws.Range('E11').Value = '  -4.19%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.60'
$ws.Range('E12').Value = '  -6.15%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.122'
$ws.Range('E13').Value = '  -7.18%  '

$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.005.06'
$ws.Range('E14').Value = '  -5.30%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.558'
$ws.Range('E15').Value = '  -5.40%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '95.85'
$ws.Range('E16').Value = '  -5.70%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.014'
$ws.Range('E17').Value = '  +0.79%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001124'
$ws.Range('E18').Value = '  -3.67%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06670'
$ws.Range('E19').Value = '  -0.78%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.98'
$ws.Range('E20').Value = '  -6.27%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.010'
$ws.Range('E21').Value = '  +0.41%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.996'
$ws.Range('E22').Value = '  -5.19%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.650.32'
$ws.Range('E23').Value = '  -2.50%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.01'
$ws.Range('E24').Value = '  -5.08%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.284'
$ws.Range('E25').Value = '  -1.33%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.57'
$ws.Range('E26').Value = '  -2.26%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.82'
$ws.Range('E27').Value = '  -5.66%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.436'
$ws.Range('E28').Value = '  -6.60%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.349'
$ws.Range('E29').Value = '  -7.81%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.52'
$ws.Range('E30').Value = '  -3.69%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.060'
$ws.Range('E31').Value = '  -7.81%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.592'
$ws.Range('E32').Value = '  -9.54%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09979'
$ws.Range('E33').Value = '  -5.53%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.880'
$ws.Range('E34').Value = '  -6.31%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.802'
$ws.Range('E35').Value = '  -3.14%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.607'
$ws.Range('E36').Value = '  -8.72%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02470'
$ws.Range('E37').Value = '  -6.32%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.330'
$ws.Range('E38').Value = '  -1.32%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06392'
$ws.Range('E39').Value = '  -6.17%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6590'
$ws.Range('E40').Value = '  -6.45%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.83'
$ws.Range('E41').Value = '  -6.03%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2083'
$ws.Range('E42').Value = '  -6.74%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.012'
$ws.Range('E43').Value = '  +0.62%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6375'
$ws.Range('E44').Value = '  -6.94%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.70'
$ws.Range('E45').Value = '  -5.83%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.217'
$ws.Range('E46').Value = '  -6.31%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.289'
$ws.Range('E47').Value = '  -4.82%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.538'
$ws.Range('E48').Value = '  -3.00%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07020'
$ws.Range('E49').Value = '  -3.05%  '

$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.138'
$ws.Range('E50').Value = '  -5.35%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.00000000321'
$ws.Range('E51').Value = '  -7.21%  '
